$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) contains values that look numeric (e.g. "1.003").
# Excel auto-converts such text to a real number when assigned via .Value,
# which would change the cell from a text cell to a numeric cell (losing the
# original trailing-zero/format-preserving text). Force the whole column to
# Text format first so the assigned strings are kept verbatim as text, then
# restore the original ("Normal") cell style so no stray number-format is left
# behind on the cells (matches the source file, where these cells carry no
# explicit style).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "28.181.47"
$ws.Range("E2").Value = "  -0.30%  "

# Row 3
$ws.Range("D3").Value = "1.868.40"
$ws.Range("E3").Value = "  +3.26%  "

# Row 4
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.28%  "

# Row 5
$ws.Range("D5").Value = "311.94"
$ws.Range("E5").Value = "  +0.06%  "

# Row 6
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.26%  "

# Row 7
$ws.Range("D7").Value = "0.5017"
$ws.Range("E7").Value = "  -2.35%  "

# Row 8
$ws.Range("D8").Value = "0.3938"
$ws.Range("E8").Value = "  -0.35%  "

# Row 9
$ws.Range("D9").Value = "0.09897"
$ws.Range("E9").Value = "  +26.59%  "

# Row 10
$ws.Range("D10").Value = "1.122"
$ws.Range("E10").Value = "  +1.17%  "

# Row 11
$ws.Range("D11").Value = "41.20"
$ws.Range("E11").Value = "  +0.56%  "

# Row 12
$ws.Range("D12").Value = "6.448"
$ws.Range("E12").Value = "  +1.38%  "

# Row 13
$ws.Range("D13").Value = "20.87"
$ws.Range("E13").Value = "  +1.86%  "

# Row 14
$ws.Range("D14").Value = "1.866.61"
$ws.Range("E14").Value = "  +3.34%  "

# Row 15
$ws.Range("D15").Value = "1.004"
$ws.Range("E15").Value = "  +0.30%  "

# Row 16
$ws.Range("D16").Value = "7.363"
$ws.Range("E16").Value = "  +0.32%  "

# Row 17
$ws.Range("D17").Value = "0.00001135"
$ws.Range("E17").Value = "  +4.97%  "

# Row 18
$ws.Range("D18").Value = "93.49"
$ws.Range("E18").Value = "  +0.90%  "

# Row 19
$ws.Range("D19").Value = "0.06649"
$ws.Range("E19").Value = "  +1.35%  "

# Row 20
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.34%  "

# Row 21
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "17.37"
$ws.Range("E21").Value = "  +0.23%  "

# Row 22
$ws.Range("D22").Value = "6.090"
$ws.Range("E22").Value = "  +1.33%  "

# Row 23
$ws.Range("D23").Value = "28.258.63"
$ws.Range("E23").Value = "  -0.11%  "

# Row 24
$ws.Range("D24").Value = "11.31"
$ws.Range("E24").Value = "  +1.50%  "

# Row 25
$ws.Range("D25").Value = "2.265"
$ws.Range("E25").Value = "  +1.79%  "

# Row 26
$ws.Range("D26").Value = "2.524"
$ws.Range("E26").Value = "  +2.63%  "

# Row 27
$ws.Range("D27").Value = "2.085.87"
$ws.Range("E27").Value = "  +3.31%  "

# Row 28
$ws.Range("D28").Value = "21.23"
$ws.Range("E28").Value = "  +3.51%  "

# Row 29
$ws.Range("D29").Value = "157.69"
$ws.Range("E29").Value = "  -1.92%  "

# Row 30
$ws.Range("D30").Value = "127.70"
$ws.Range("E30").Value = "  -0.09%  "

# Row 31
$ws.Range("D31").Value = "0.1058"
$ws.Range("E31").Value = "  -3.62%  "

# Row 32
$ws.Range("D32").Value = "1.051"
$ws.Range("E32").Value = "  -0.92%  "

# Row 33
$ws.Range("D33").Value = "5.613"
$ws.Range("E33").Value = "  +0.85%  "

# Row 34
$ws.Range("D34").Value = "3.610"
$ws.Range("E34").Value = "  -1.10%  "

# Row 35
$ws.Range("D35").Value = "0.06784"
$ws.Range("E35").Value = "  -5.31%  "

# Row 36
$ws.Range("D36").Value = "9.397"
$ws.Range("E36").Value = "  +1.94%  "

# Row 37
$ws.Range("D37").Value = "0.02385"
$ws.Range("E37").Value = "  +1.36%  "

# Row 38
$ws.Range("D38").Value = "0.2179"
$ws.Range("E38").Value = "  +0.20%  "

# Row 39
$ws.Range("D39").Value = "4.997"
$ws.Range("E39").Value = "  -0.58%  "

# Row 40
$ws.Range("D40").Value = "11.46"
$ws.Range("E40").Value = "  -1.04%  "

# Row 41
$ws.Range("D41").Value = "0.6279"
$ws.Range("E41").Value = "  +1.55%  "

# Row 42
$ws.Range("D42").Value = "1.170"
$ws.Range("E42").Value = "  +0.93%  "

# Row 43
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  +0.32%  "

# Row 44
$ws.Range("D44").Value = "13.52"
$ws.Range("E44").Value = "  +1.93%  "

# Row 45
$ws.Range("D45").Value = "0.5996"
$ws.Range("E45").Value = "  +0.20%  "

# Row 46
$ws.Range("D46").Value = "1.278"
$ws.Range("E46").Value = "  -2.24%  "

# Row 47
$ws.Range("D47").Value = "3.664"
$ws.Range("E47").Value = "  -1.99%  "

# Row 48
$ws.Range("D48").Value = "124.63"
$ws.Range("E48").Value = "  -0.38%  "

# Row 49
$ws.Range("D49").Value = "1.985"
$ws.Range("E49").Value = "  +3.42%  "

# Row 50
$ws.Range("D50").Value = "1.195"
$ws.Range("E50").Value = "  -1.14%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.06845"
$ws.Range("E51").Value = "  +0.71%  "

# Restore the default cell style on the Price column so it matches the
# original (un-styled) cells now that the text values are locked in.
$priceRange.Style = "Normal"
